# Adding code for Adaptive algorithm
# Append the next Latency/Runtime data point (50 ms -> 1.531 sec) as row 8
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 50
$ws.Range("B8").Value = 1.531

# Match Excel's post-edit active selection (cursor parked on the newly
# entered cell, as it would be after typing the value and hitting Enter)
$ws.Range("B8").Select()
